# Update crypto price/volume figures in the cryptos list (GitHub Actions refresh).
# Column D ("Price") holds values that look numeric but must stay plain text
# (the sheet stores them as inline strings, e.g. "42.342.94"), so we force the
# cell to Text format before writing - otherwise Excel would silently coerce a
# value like "230.63" into a real number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.342.94'
$ws.Range("E2").Value = '  +1.91%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.218.57'
$ws.Range("E3").Value = '  -0.15%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.63'
$ws.Range("E5").Value = '  +0.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("E6").Value = '  -0.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.69'
$ws.Range("E7").Value = '  -0.78%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.402'
$ws.Range("E9").Value = '  -0.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0903'
$ws.Range("E10").Value = '  +2.41%  '

$ws.Range("E11").Value = '  +0.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.545.39'
$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.49'
$ws.Range("E13").Value = '  -0.56%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.14'
$ws.Range("E14").Value = '  +2.87%  '

$ws.Range("E15").Value = '  +1.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.796'
$ws.Range("E16").Value = '  +0.15%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.238.31'
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.161.97'
$ws.Range("E18").Value = '  +1.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0937'
$ws.Range("E19").Value = '  +5.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.16'
$ws.Range("E20").Value = '  +2.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.07'
$ws.Range("E21").Value = '  -0.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '243.89'
$ws.Range("E22").Value = '  -2.06%  '

$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.39'
$ws.Range("E24").Value = '  +1.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.29'
$ws.Range("E25").Value = '  -0.77%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.65'
$ws.Range("E26").Value = '  +0.99%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.54'
$ws.Range("E27").Value = '  +1.06%  '

$ws.Range("E28").Value = '  -1.14%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.33'
$ws.Range("E29").Value = '  +2.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.44'
$ws.Range("E30").Value = '  +1.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.66'
$ws.Range("E31").Value = '  +4.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.121'
$ws.Range("E32").Value = '  -1.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.98'
$ws.Range("E33").Value = '  +0.12%  '

$ws.Range("E34").Value = '  +0.12%  '

$ws.Range("E35").Value = '  +4.71%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.34'
$ws.Range("E36").Value = '  -2.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.54'
$ws.Range("E37").Value = '  -3.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.34'
$ws.Range("E38").Value = '  -1.24%  '

$ws.Range("E39").Value = '  +6.50%  '

$ws.Range("E40").Value = '  -0.14%  '

$ws.Range("E41").Value = '  -2.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.55'
$ws.Range("E42").Value = '  -1.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0958'
$ws.Range("E43").Value = '  -2.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.20'
$ws.Range("E44").Value = '  +0.76%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.79'
$ws.Range("E45").Value = '  -2.19%  '

$ws.Range("E46").Value = '  -9.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.456.51'
$ws.Range("E47").Value = '  -0.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.76'
$ws.Range("E48").Value = '  -0.98%  '

$ws.Range("E49").Value = '  -0.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '15.95'
$ws.Range("E50").Value = '  -2.69%  '

$ws.Range("E51").Value = '  +3.16%  '
